$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (Great Sitkin): stop date and duration increased by one day
$ws.Range("C13").Value2 = 45492
$ws.Range("D13").Value2 = 1151

# Remove four eruption records that lacked confirmed stop dates
# (delete from the bottom up so row numbers of earlier rows stay valid)
$ws.Rows.Item(76).Delete()
$ws.Rows.Item(75).Delete()
$ws.Rows.Item(68).Delete()
$ws.Rows.Item(46).Delete()
